# Auto-generated update of LeveProfits data cells (H..N) across all sheets
# per the diff: value-only corrections, no formula/format changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2242.4375
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 2325.2666
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 2325.2666
$ws.Range("M32").Value = -674
$ws.Range("N32").Value = -2977.2666

$ws.Range("H62").Value = 4431
$ws.Range("I62").Value = 3895
$ws.Range("K62").Value = 3895
$ws.Range("M62").Value = -3271

$ws.Range("H64").Value = 4136
$ws.Range("I64").Value = 4248
$ws.Range("K64").Value = 4248
$ws.Range("M64").Value = -4000

$ws.Range("H65").Value = 4431
$ws.Range("I65").Value = 3895
$ws.Range("K65").Value = 19475
$ws.Range("M65").Value = -16355

$ws.Range("H67").Value = 4136
$ws.Range("I67").Value = 4248
$ws.Range("K67").Value = 4248
$ws.Range("M67").Value = -3390

$ws.Range("H112").Value = 3549.25
$ws.Range("I112").Value = 1025
$ws.Range("J112").Value = 4180.3125
$ws.Range("K112").Value = 3075
$ws.Range("L112").Value = 12540.9375
$ws.Range("M112").Value = -1967
$ws.Range("N112").Value = -14756.9375

$ws.Range("H121").Value = 913.94116
$ws.Range("I121").Value = 1500
$ws.Range("J121").Value = 877.3125
$ws.Range("K121").Value = 4500
$ws.Range("L121").Value = 2631.9375
$ws.Range("M121").Value = -2753
$ws.Range("N121").Value = -6125.9375

$ws.Range("H131").Value = 863.5714
$ws.Range("I131").Value = 863.5714
$ws.Range("K131").Value = 2590.7142
$ws.Range("M131").Value = 2449.2858

$ws.Range("H132").Value = 6949403.5
$ws.Range("I132").Value = 10105230
$ws.Range("J132").Value = 6586.933
$ws.Range("K132").Value = 30315690
$ws.Range("L132").Value = 19760.799
$ws.Range("M132").Value = -30313160
$ws.Range("N132").Value = -24820.799

$ws.Range("H137").Value = 1382.4
$ws.Range("J137").Value = 1739.5
$ws.Range("L137").Value = 5218.5
$ws.Range("N137").Value = -10318.5

$ws.Range("H138").Value = 1548.8282
$ws.Range("I138").Value = 935.125
$ws.Range("J138").Value = 1667.1326
$ws.Range("K138").Value = 2805.375
$ws.Range("L138").Value = 5001.3978
$ws.Range("M138").Value = 2334.625
$ws.Range("N138").Value = -15281.3978

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M22").ClearContents()
$ws.Range("H22").Value = 109
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 109
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 109
$ws.Range("N22").Value = -707

$ws.Range("H32").Value = 3416.744
$ws.Range("I32").Value = 3131.8918
$ws.Range("K32").Value = 3131.8918
$ws.Range("M32").Value = -2844.8918

$ws.Range("H41").Value = 3286.8
$ws.Range("I41").Value = 3286.8
$ws.Range("K41").Value = 3286.8
$ws.Range("M41").Value = -2872.8

$ws.Range("H74").Value = 1491.3784
$ws.Range("I74").Value = 854.75
$ws.Range("J74").Value = 2240.353
$ws.Range("K74").Value = 854.75
$ws.Range("L74").Value = 2240.353
$ws.Range("M74").Value = 19.25
$ws.Range("N74").Value = -3988.353

$ws.Range("H77").Value = 1491.3784
$ws.Range("I77").Value = 854.75
$ws.Range("J77").Value = 2240.353
$ws.Range("K77").Value = 4273.75
$ws.Range("L77").Value = 11201.765
$ws.Range("M77").Value = 94.25
$ws.Range("N77").Value = -19937.765

$ws.Range("H132").Value = 1554.8223
$ws.Range("I132").Value = 1269.9
$ws.Range("J132").Value = 2124.6667
$ws.Range("K132").Value = 3809.7
$ws.Range("L132").Value = 6374.000100000001
$ws.Range("M132").Value = -1279.7
$ws.Range("N132").Value = -11434.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 11364327
$ws.Range("I94").Value = 14706494
$ws.Range("J94").Value = 956.8
$ws.Range("K94").Value = 14706494
$ws.Range("L94").Value = 956.8
$ws.Range("M94").Value = -14706043
$ws.Range("N94").Value = -1858.8

$ws.Range("H134").Value = 6247.1924
$ws.Range("I134").Value = 1420.65
$ws.Range("J134").Value = 22335.666
$ws.Range("K134").Value = 4261.950000000001
$ws.Range("L134").Value = 67006.99800000001
$ws.Range("M134").Value = -1726.950000000001
$ws.Range("N134").Value = -72076.99800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1240
$ws.Range("I31").Value = 982.2
$ws.Range("J31").Value = 1884.5
$ws.Range("K31").Value = 982.2
$ws.Range("L31").Value = 1884.5
$ws.Range("M31").Value = -687.2
$ws.Range("N31").Value = -2474.5

$ws.Range("H34").Value = 1240
$ws.Range("I34").Value = 982.2
$ws.Range("J34").Value = 1884.5
$ws.Range("K34").Value = 982.2
$ws.Range("L34").Value = 1884.5
$ws.Range("M34").Value = -780.2
$ws.Range("N34").Value = -2288.5

$ws.Range("H58").Value = 973.5
$ws.Range("I58").Value = 875
$ws.Range("J58").Value = 1308.4
$ws.Range("K58").Value = 875
$ws.Range("L58").Value = 1308.4
$ws.Range("M58").Value = -672
$ws.Range("N58").Value = -1714.4

$ws.Range("H62").Value = 66669332
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376

$ws.Range("H65").Value = 66669332
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880

$ws.Range("H86").Value = 3718308.5
$ws.Range("J86").Value = 26162.445
$ws.Range("L86").Value = 26162.445
$ws.Range("N86").Value = -28408.445

$ws.Range("H89").Value = 3718308.5
$ws.Range("J89").Value = 26162.445
$ws.Range("L89").Value = 130812.225
$ws.Range("N89").Value = -142044.225

$ws.Range("H132").Value = 1416.1842
$ws.Range("I132").Value = 945.7742
$ws.Range("J132").Value = 3499.4285
$ws.Range("K132").Value = 2837.3226
$ws.Range("L132").Value = 10498.2855
$ws.Range("M132").Value = -307.3226
$ws.Range("N132").Value = -15558.2855

$ws.Range("H134").Value = 1088.7333
$ws.Range("I134").Value = 937
$ws.Range("K134").Value = 2811
$ws.Range("M134").Value = -276

$ws.Range("H136").Value = 973.5
$ws.Range("I136").Value = 875
$ws.Range("J136").Value = 1308.4
$ws.Range("K136").Value = 2625
$ws.Range("L136").Value = 3925.2
$ws.Range("M136").Value = -75
$ws.Range("N136").Value = -9025.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 20836282
$ws.Range("I131").Value = 83333540
$ws.Range("J131").Value = 3866.5
$ws.Range("K131").Value = 250000620
$ws.Range("L131").Value = 11599.5
$ws.Range("M131").Value = -249995580
$ws.Range("N131").Value = -21679.5

$ws.Range("H139").Value = 1898.1482
$ws.Range("I139").Value = 2034.875
$ws.Range("J139").Value = 1699.2727
$ws.Range("K139").Value = 6104.625
$ws.Range("L139").Value = 5097.8181
$ws.Range("M139").Value = -964.625
$ws.Range("N139").Value = -15377.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1687.3721
$ws.Range("I132").Value = 1148.4286
$ws.Range("J132").Value = 2693.4
$ws.Range("K132").Value = 3445.2858
$ws.Range("L132").Value = 8080.200000000001
$ws.Range("M132").Value = -915.2857999999997
$ws.Range("N132").Value = -13140.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 5037
$ws.Range("J41").Value = 5037
$ws.Range("L41").Value = 5037
$ws.Range("N41").Value = -5913

$ws.Range("H132").Value = 17903.541
$ws.Range("I132").Value = 987.53845
$ws.Range("K132").Value = 2962.61535
$ws.Range("M132").Value = -432.61535

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2268.0952
$ws.Range("I132").Value = 2283.2974
$ws.Range("K132").Value = 6849.8922
$ws.Range("M132").Value = -4319.8922

$ws.Range("H136").Value = 469.9
$ws.Range("I136").Value = 205.76471
$ws.Range("J136").Value = 1966.6666
$ws.Range("K136").Value = 617.29413
$ws.Range("L136").Value = 5899.9998
$ws.Range("M136").Value = 1932.70587
$ws.Range("N136").Value = -10999.9998

$ws.Range("H141").Value = 25059.584
$ws.Range("J141").Value = 25059.584
$ws.Range("L141").Value = 25059.584
$ws.Range("N141").Value = -35419.584
